# Add the new "CAZyme domains" worksheet (ANOVA results for CAZyme domain substrates),
# following the same layout/formatting used by the existing "litterChemistry" sheet.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("litterChemistry")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CAZyme domains"

# Copy formatting (header/border styles, column widths, etc.) from litterChemistry.
# litterChemistry has 11 rows (1 header + 10 data); this new sheet needs 15 rows
# (1 header + 14 data), so stamp the header format once and the data-row format
# across all of rows 2-15.
$src.Range("A1:H11").Copy()
$ws.Range("A1:H11").PasteSpecial(-4122) # xlPasteFormats
$src.Range("A2:H2").Copy()
$ws.Range("A12:H15").PasteSpecial(-4122) # xlPasteFormats

# Match the (auto-fit) column widths used on litterChemistry's matching columns.
for ($i = 1; $i -le 8; $i++) {
    $ws.Cells.Item(1, $i).ColumnWidth = $src.Cells.Item(1, $i).ColumnWidth
}

# --- Fill in substrate names (column A) in the same order they were authored ---
$ws.Range("A2").Value = "Hemicellulose"
$ws.Range("A3").Value = "Lignin"
$ws.Range("A4").Value = "Polysaccharide"
$ws.Range("A5").Value = "Oligosaccharides"
$ws.Range("A6").Value = "Cell_wall"
$ws.Range("A7").Value = "Inulin"
$ws.Range("A8").Value = "Starch"
$ws.Range("A9").Value = "Trehalose"
$ws.Range("A10").Value = "Cellulose"
$ws.Range("A11").Value = "Pectin"
$ws.Range("A12").Value = "Glycogen"
$ws.Range("A15").Value = "Total"
$ws.Range("A13").Value = "Peptidoglycan"
$ws.Range("A14").Value = "Chitin"
$ws.Range("A1").Value = "Substrate"

# --- Header row ---
$ws.Range("B1").Value = "timePoint"
$ws.Range("C1").Value = "Vegetation"
$ws.Range("D1").Value = "Precipitation"
$ws.Range("E1").Value = "timePoint x Precipitation"
$ws.Range("F1").Value = "timePoint x Vegetation"
$ws.Range("G1").Value = "Vegetation x Precipitation"
$ws.Range("H1").Value = "Three-way"

# --- ANOVA significance results ---
$ws.Range("B2").Value = "***"
$ws.Range("C2").Value = "***"
$ws.Range("D2").Value = "o"
$ws.Range("E2").Value = "o"
$ws.Range("F2").Value = "*"
$ws.Range("G2").Value = "o"
$ws.Range("H2").Value = "o"

$ws.Range("B3").Value = "*"
$ws.Range("C3").Value = "***"
$ws.Range("D3").Value = "***"
$ws.Range("E3").Value = "**"
$ws.Range("F3").Value = "o"
$ws.Range("G3").Value = "o"
$ws.Range("H3").Value = "*"

$ws.Range("B4").Value = "***"
$ws.Range("C4").Value = "***"
$ws.Range("D4").Value = "o"
$ws.Range("E4").Value = "o"
$ws.Range("F4").Value = "o"
$ws.Range("G4").Value = "o"
$ws.Range("H4").Value = "*"

$ws.Range("B5").Value = "*"
$ws.Range("C5").Value = "***"
$ws.Range("D5").Value = "o"
$ws.Range("E5").Value = "o"
$ws.Range("F5").Value = "***"
$ws.Range("G5").Value = "o"
$ws.Range("H5").Value = "o"

$ws.Range("B6").Value = "***"
$ws.Range("C6").Value = "***"
$ws.Range("D6").Value = "o"
$ws.Range("E6").Value = "o"
$ws.Range("F6").Value = "o"
$ws.Range("G6").Value = "o"
$ws.Range("H6").Value = "*"

$ws.Range("B7").Value = "***"
$ws.Range("C7").Value = "*"
$ws.Range("D7").Value = "o"
$ws.Range("E7").Value = "o"
$ws.Range("F7").Value = "o"
$ws.Range("G7").Value = "***"
$ws.Range("H7").Value = "**"

$ws.Range("B8").Value = "o"
$ws.Range("C8").Value = "o"
$ws.Range("D8").Value = "o"
$ws.Range("E8").Value = "o"
$ws.Range("F8").Value = "o"
$ws.Range("G8").Value = "*"
$ws.Range("H8").Value = "o"

$ws.Range("B9").Value = "***"
$ws.Range("C9").Value = "**"
$ws.Range("D9").Value = "o"
$ws.Range("E9").Value = "o"
$ws.Range("F9").Value = "o"
$ws.Range("G9").Value = "o"
$ws.Range("H9").Value = "o"

$ws.Range("B10").Value = "***"
$ws.Range("C10").Value = "**"
$ws.Range("D10").Value = "o"
$ws.Range("E10").Value = "o"
$ws.Range("F10").Value = "o"
$ws.Range("G10").Value = "***"
$ws.Range("H10").Value = "***"

$ws.Range("B11").Value = "**"
$ws.Range("C11").Value = "o"
$ws.Range("D11").Value = "**"
$ws.Range("E11").Value = "o"
$ws.Range("F11").Value = "o"
$ws.Range("G11").Value = "***"
$ws.Range("H11").Value = "o"

$ws.Range("B12").Value = "**"
$ws.Range("C12").Value = "o"
$ws.Range("D12").Value = "o"
$ws.Range("E12").Value = "o"
$ws.Range("F12").Value = "**"
$ws.Range("G12").Value = "o"
$ws.Range("H12").Value = "o"

$ws.Range("B13").Value = "**"
$ws.Range("C13").Value = "o"
$ws.Range("D13").Value = "o"
$ws.Range("E13").Value = "o"
$ws.Range("F13").Value = "o"
$ws.Range("G13").Value = "o"
$ws.Range("H13").Value = "o"

$ws.Range("B14").Value = "**"
$ws.Range("C14").Value = "o"
$ws.Range("D14").Value = "o"
$ws.Range("E14").Value = "o"
$ws.Range("F14").Value = "o"
$ws.Range("G14").Value = "o"
$ws.Range("H14").Value = "o"

$ws.Range("B15").Value = "**"
$ws.Range("C15").Value = "***"
$ws.Range("D15").Value = "o"
$ws.Range("E15").Value = "o"
$ws.Range("F15").Value = "*"
$ws.Range("G15").Value = "o"
$ws.Range("H15").Value = "*"

# --- View/selection state to match the saved workbook ---
# litterChemistry is no longer the active/selected sheet; its lingering selection
# is reset to the header row.
$src.Activate()
$src.Range("A1:H1").Select()

# The new sheet becomes the active tab, with the last selected cell at B16.
$ws.Activate()
$ws.Range("B16").Select()

Write-Host "CAZyme domains sheet added"
